# Updated cryptos list on Fri May 17 11:45:04 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for
# each coin row to their latest scraped values, and swaps the Filecoin /
# Mantle rows (36 and 37), which changed rank order, carrying their new
# prices along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text so Excel does not silently coerce
# numeric-looking strings (e.g. "578.70") into floating point numbers,
# which would both lose the trailing zero and change the cell type.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "66.328.38"
$ws.Range("E2").Value = "  +0.28%  "
Set-TextValue "D3" "3.033.32"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue "D5" "578.70"
$ws.Range("E5").Value = "  -0.44%  "
Set-TextValue "D6" "168.29"
$ws.Range("E6").Value = "  +3.32%  "
Set-TextValue "D8" "3.031.98"
$ws.Range("E8").Value = "  +1.27%  "
Set-TextValue "D9" "0.522"
$ws.Range("E9").Value = "  +0.68%  "
Set-TextValue "D10" "6.67"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("E11").Value = "  -1.17%  "
Set-TextValue "D12" "0.488"
$ws.Range("E12").Value = "  +7.25%  "
$ws.Range("E13").Value = "  -1.16%  "
Set-TextValue "D14" "36.74"
$ws.Range("E14").Value = "  +6.44%  "
$ws.Range("E15").Value = "  +0.15%  "
Set-TextValue "D16" "66.312.05"
$ws.Range("E16").Value = "  +0.24%  "
Set-TextValue "D17" "3.535.37"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("E18").Value = "  +4.70%  "
Set-TextValue "D19" "16.55"
$ws.Range("E19").Value = "  +19.65%  "
Set-TextValue "D20" "3.031.68"
$ws.Range("E20").Value = "  +0.99%  "
Set-TextValue "D21" "467.24"
$ws.Range("E21").Value = "  +3.30%  "
Set-TextValue "D22" "0.712"
$ws.Range("E22").Value = "  +4.13%  "
$ws.Range("E23").Value = "  +0.95%  "
Set-TextValue "D24" "83.10"
$ws.Range("E24").Value = "  +0.95%  "
Set-TextValue "D25" "12.74"
$ws.Range("E25").Value = "  +3.83%  "
$ws.Range("E26").Value = "  -1.18%  "
Set-TextValue "D27" "10.04"
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("E28").Value = "  +0.03%  "
Set-TextValue "D29" "8.20"
$ws.Range("E29").Value = "  +1.15%  "
Set-TextValue "D30" "2.43"
$ws.Range("E30").Value = "  +1.16%  "
Set-TextValue "D31" "2.66"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("E32").Value = "  +6.86%  "
Set-TextValue "D33" "0.0₃0997"
$ws.Range("E33").Value = "  -3.46%  "
Set-TextValue "D34" "28.12"
$ws.Range("E34").Value = "  +3.42%  "
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  -0.13%  "

# Rows 36 and 37 swap content: Filecoin (previously row 36) moves to row
# 37, and Mantle (previously row 37) moves to row 36, each carrying
# refreshed price/volume figures.
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D36" "0.994"
$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D37" "5.85"
$ws.Range("E37").Value = "  +0.76%  "

Set-TextValue "D38" "48.27"
$ws.Range("E38").Value = "  +9.71%  "
$ws.Range("E39").Value = "  -0.27%  "
Set-TextValue "D40" "49.56"
$ws.Range("E40").Value = "  +0.07%  "
Set-TextValue "D41" "0.313"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("E43").Value = "  +2.56%  "
Set-TextValue "D44" "2.85"
$ws.Range("E44").Value = "  -3.53%  "
$ws.Range("E45").Value = "  +0.74%  "
Set-TextValue "D46" "379.35"
$ws.Range("E46").Value = "  -5.53%  "
Set-TextValue "D47" "2.705.41"
$ws.Range("E47").Value = "  -2.23%  "
Set-TextValue "D48" "134.22"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("E49").Value = "  +0.05%  "
Set-TextValue "D50" "24.47"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("E51").Value = "  +4.17%  "
